# 440-RBI-EI-DB-SAR-REC-NON-RNI-CTRFD-DL-MD-TR-1-EarlyRePayment-Loanproduct4.xlsx
# "code refactoring and loan accounting and charges added"

$wb  = $excel.ActiveWorkbook
$wsIn  = $wb.Worksheets.Item("ProductLoan_Input")
$wsOut = $wb.Worksheets.Item("ProductLoan_Output")

# ---------------------------------------------------------------------------
# 1. Product name renamed (drives both the Input and Output sheet headers,
#    which mirror the same product name value in B1).
# ---------------------------------------------------------------------------
$productName = "440-RBI-EI-DB-SAR-REC-NON-RNI-CTRFD-DL-MD-TR-1-EarlyRePayment"
$wsIn.Range("B1").Value  = $productName
$wsOut.Range("B1").Value = $productName

# ---------------------------------------------------------------------------
# 2. Short name changed from text "kar7" to the numeric code 440.
# ---------------------------------------------------------------------------
$wsIn.Range("B3").Value = 440

# ---------------------------------------------------------------------------
# 3. "nominalinterestratedefault" value dropped from 12 down to 1.
# ---------------------------------------------------------------------------
$wsIn.Range("B11").Value = 1

# ---------------------------------------------------------------------------
# 4. New loan-accounting / charge mapping rows appended (31-42), copying the
#    existing alternating row style (label col = style used by A10,
#    value col = style used by B10) down across the new block first, then
#    filling in the text.
# ---------------------------------------------------------------------------
$wsIn.Range("A10:B10").Copy()
$wsIn.Range("A31:B42").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$accountingRows = @(
  @("fundsource",                     "Cash"),
  @("loanprotfolio",                  "Loan portfolio "),
  @("interestreceivable",             "Interest Receivable "),
  @("penaltiesreceivable",            "Penalties Receivable "),
  @("transferinsuspense",             "Transfer in Suspence "),
  @("feesreceivable",                 "Fees Receivable"),
  @("incomefrominterest",             "Income from interest"),
  @("incomefrompenalties",            "Income from penalties"),
  @("incomefromfees",                 "Income from fees"),
  @("incomefromrecoveryrepayments",   "Income from recovery repayments"),
  @("loseswrittenoff",                "Losses Writtenoff "),
  @("overpaymentliability",           "Overpayment Liability")
)

$row = 31
foreach ($pair in $accountingRows) {
  $wsIn.Cells.Item($row, 1).Value = $pair[0]
  $wsIn.Cells.Item($row, 2).Value = $pair[1]
  $row = $row + 1
}

# ---------------------------------------------------------------------------
# 5. View state: selection moves down to the new rows on the input sheet,
#    and the output sheet becomes the active/selected tab.
# ---------------------------------------------------------------------------
[void]$wsIn.Range("A40").Select()
[void]$wsOut.Range("B1").Select()
[void]$wsOut.Activate()
